$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time number format (style index 1 in the original file) from
# the last existing row (A115) onto the two new date cells, then overwrite
# the values. Using Copy() (format-only carries over, value is overwritten
# right after) reuses the existing style record instead of minting a new one.
$ws.Range("A115").Copy($ws.Range("A116"))
$ws.Range("A115").Copy($ws.Range("A117"))

# --- Row 116 ---
$ws.Range("A116").Value = 45489.2916666667
$ws.Range("B116").Value = 0
$ws.Range("C116").Value = 0.670000016689301
$ws.Range("D116").Value = 0.670000016689301
$ws.Range("E116").Value = 0.670000016689301
$ws.Range("F116").Value = 0.670000016689301

# G/H hold the values as text (shared strings in the saved file). Marking the
# cell as Text before assigning keeps Excel from re-parsing the numeric-
# looking string back into a number; resetting to Normal afterwards restores
# the default (unstyled) cell appearance the source file uses.
$ws.Range("G116").NumberFormat = "@"
$ws.Range("G116").Value = "0.670000016689301"
$ws.Range("G116").Style = "Normal"

$ws.Range("H116").NumberFormat = "@"
$ws.Range("H116").Value = "BWZ.MI"
$ws.Range("H116").Style = "Normal"

# --- Row 117 ---
$ws.Range("A117").Value = 45490.6496180556
$ws.Range("B117").Value = 10551
$ws.Range("C117").Value = 0.709999978542328
$ws.Range("D117").Value = 0.639999985694885
$ws.Range("E117").Value = 0.670000016689301
$ws.Range("F117").Value = 0.680000007152557

$ws.Range("G117").NumberFormat = "@"
$ws.Range("G117").Value = "0.680000007152557"
$ws.Range("G117").Style = "Normal"

$ws.Range("H117").NumberFormat = "@"
$ws.Range("H117").Value = "BWZ.MI"
$ws.Range("H117").Style = "Normal"
